$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is "Drag/drop-funktion" -> "Skapa funktion som returnerar en array med 4 slumptal 0-7"
# Mark it as finished ("Färdig") and record 1 hour of actual time.
$ws.Range("C19").Value = "Färdig"
$ws.Range("E19").Value = 1

# Update the active selection to match the new state
$ws.Range("E24").Select()
